$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) and "全部类型" (sheet4) both contain the same
# four rows (2-5) with "想去人数" (want-to-go count) values in column F
# that need to be bumped to their latest scraped totals.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 159
    $ws.Range("F3").Value = 61
    $ws.Range("F4").Value = 264
    $ws.Range("F5").Value = 4024
}
